$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must be forced to stay as text
# (otherwise Excel auto-converts "573.93" -> number 573.93, losing the
# original inline-string type and any trailing zeros, e.g. "24.10" -> 24.1)
$forceTextCells = @("D5", "D6", "D7", "D14", "D18", "D19", "D20", "D22", "D24", "D26", "D28", "D30", "D32", "D33", "D35", "D36", "D37", "D41", "D43", "D44", "D45", "D46", "D47", "D50", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.115.45'
$ws.Range("E2").Value = '  -4.01%  '

$ws.Range("D3").Value = '3.319.33'
$ws.Range("E3").Value = '  -0.53%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '573.93'
$ws.Range("E5").Value = '  -2.58%  '

$ws.Range("D6").Value = '181.43'
$ws.Range("E6").Value = '  -4.11%  '

$ws.Range("D7").Value = '0.616'
$ws.Range("E7").Value = '  +2.17%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -3.25%  '

$ws.Range("E10").Value = '  -1.34%  '

$ws.Range("E11").Value = '  -2.86%  '

$ws.Range("D12").Value = '3.901.15'
$ws.Range("E12").Value = '  -0.41%  '

$ws.Range("E13").Value = '  -1.37%  '

$ws.Range("D14").Value = '26.68'

$ws.Range("D15").Value = '66.258.94'
$ws.Range("E15").Value = '  -3.87%  '

$ws.Range("E16").Value = '  -2.05%  '

$ws.Range("D17").Value = '3.317.36'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").Value = '436.51'
$ws.Range("E18").Value = '  -2.85%  '

$ws.Range("D19").Value = '5.67'

$ws.Range("D20").Value = '13.55'
$ws.Range("E20").Value = '  -1.65%  '

$ws.Range("E21").Value = '  -3.01%  '

$ws.Range("D22").Value = '73.19'
$ws.Range("E22").Value = '  -3.47%  '

$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").Value = '0.520'
$ws.Range("E24").Value = '  -0.54%  '

$ws.Range("E25").Value = '  -3.94%  '

$ws.Range("D26").Value = '0.192'
$ws.Range("E26").Value = '  +0.81%  '

$ws.Range("E27").Value = '  -3.45%  '

$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.22%  '

$ws.Range("E29").Value = '  -3.01%  '

$ws.Range("D30").Value = '22.73'
$ws.Range("E30").Value = '  -2.16%  '

$ws.Range("E31").Value = '  +0.03%  '

$ws.Range("D32").Value = '5.25'
$ws.Range("E32").Value = '  -4.06%  '

$ws.Range("D33").Value = '6.76'
$ws.Range("E33").Value = '  -2.48%  '

$ws.Range("E34").Value = '  -4.42%  '

$ws.Range("D35").Value = '160.50'
$ws.Range("E35").Value = '  -1.73%  '

$ws.Range("D36").Value = '1.48'
$ws.Range("E36").Value = '  -4.30%  '

$ws.Range("D37").Value = '27.78'
$ws.Range("E37").Value = '  +2.44%  '

$ws.Range("E38").Value = '  -7.11%  '

$ws.Range("D39").Value = '2.837.67'
$ws.Range("E39").Value = '  +5.33%  '

$ws.Range("E40").Value = '  -0.57%  '

$ws.Range("D41").Value = '4.43'
$ws.Range("E41").Value = '  -3.82%  '

$ws.Range("E42").Value = '  -4.47%  '

$ws.Range("D43").Value = '40.20'
$ws.Range("E43").Value = '  -2.31%  '

$ws.Range("D44").Value = '0.0666'
$ws.Range("E44").Value = '  -2.57%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '24.10'
$ws.Range("E45").Value = '  -3.79%  '

$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '2.34'
$ws.Range("E46").Value = '  -5.56%  '

$ws.Range("D47").Value = '323.46'
$ws.Range("E47").Value = '  -2.34%  '

$ws.Range("E48").Value = '  -3.50%  '

$ws.Range("E49").Value = '  +1.26%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '6.15'
$ws.Range("E50").Value = '  -2.51%  '

$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").Value = '0.972'
$ws.Range("E51").Value = '  -2.93%  '

# Restore default (unstyled) style on cells we temporarily reformatted as text,
# so the text-forcing trick leaves no visible style/format change behind.
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
